$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.282712168918807
$ws.Range("D2").Value = 2.546408760603863
$ws.Range("E2").Value = 0.8172749181555058
$ws.Range("B3").Value = 0.8590649634300813
$ws.Range("C3").Value = 0.9162403238120435
$ws.Range("D3").Value = 1.985009589357345
$ws.Range("E3").Value = 0.8590649634300818
$ws.Range("B4").Value = 0.7693344741372451
$ws.Range("C4").Value = 1.263827901423533
$ws.Range("D4").Value = 2.850324912663027
$ws.Range("E4").Value = 0.7751513876517242
$ws.Range("B5").Value = 0.4113856295130186
$ws.Range("D5").Value = 1.501212442401527
$ws.Range("E5").Value = 0.5431842522554702
$ws.Range("B6").Value = 0.6632026352630535
$ws.Range("D6").Value = 1.105282346030013
$ws.Range("E6").Value = 0.6632026352630529
$ws.Range("B7").Value = 0.36866429029728
$ws.Range("C7").Value = 0.1994065337379439
$ws.Range("D7").Value = 1.091702400591007
$ws.Range("E7").Value = 0.672249785793101
$ws.Range("D8").Value = 4.38379157803626
$ws.Range("E8").Value = 0.7197748816916139
$ws.Range("D9").Value = 3.198181802211544
$ws.Range("E9").Value = 0.7930384151307969
$ws.Range("B10").Value = 0.6187056804012239
$ws.Range("C10").Value = 1.786984753496657
$ws.Range("D10").Value = 3.921996075582012
$ws.Range("E10").Value = 0.712744398541705
$ws.Range("B11").Value = 0.419372129527483
$ws.Range("C11").Value = 71.33157886927523
$ws.Range("E11").Value = 0.5381202272564958
$ws.Range("B12").Value = 0.6524561566930798
$ws.Range("C12").Value = 55.18711579451868
$ws.Range("D12").Value = 176.2549992840693
$ws.Range("E12").Value = 0.6524561566930794
$ws.Range("B13").Value = 0.4994516281317321
$ws.Range("C13").Value = 71.42942567176385
$ws.Range("D13").Value = 128.6531648416061
$ws.Range("E13").Value = 0.7635456769531122
$ws.Range("B14").Value = 0.5794894422823409
$ws.Range("C14").Value = 0.2987436522198714
$ws.Range("D14").Value = 0.3599677653554678
$ws.Range("E14").Value = 0.6492545193630047
$ws.Range("D15").Value = 0.2972480973687998
$ws.Range("E15").Value = 0.7157323746001714
$ws.Range("B16").Value = 0.3898741030172687
$ws.Range("C16").Value = 0.5304362184867907
$ws.Range("D16").Value = 0.3108195518084244
$ws.Range("E16").Value = 0.9248938204621666
